# "update rencana bon voyage" - swap the cadet on this PDH uniform-sizing
# sheet (the MERGEFIELD result text cached in the document body) from
# SANDY SATRIA WIDJAYA to ARYA NABIL ABYAN, along with his NO / NIM /
# class / uniform-size figures. The sheet repeats the same merge block
# twice, so every replacement below hits two occurrences.

$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

function Replace-WholeWord($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# MERGEFIELD NAMA
Replace-All "SANDY SATRIA WIDJAYA" "ARYA NABIL ABYAN"
# MERGEFIELD TOPI (NIM)
Replace-All "2020.01.2.0008" "2020.01.1.0024"
# MERGEFIELD KELAS
Replace-All "UHT-TEKNIKA" "UHT-NAUTIKA"

# Short numeric fields need whole-word matching so we don't clobber
# substrings inside other numbers/dates in the sheet.
# MERGEFIELD NO
Replace-WholeWord "2" "24"
# MERGEFIELD UB_1
Replace-WholeWord "43" "40"
# MERGEFIELD UB_3
Replace-WholeWord "15" "16"
# MERGEFIELD UB_7
Replace-WholeWord "70" "68"
# MERGEFIELD UB_8
Replace-WholeWord "36" "37"
